# edit.ps1
# Applies the "Project 2 Complete... again" revision to Project Plan.docx
#
# Summary of changes (top to bottom of the document):
#  1. Insert a new paragraph "Load results graph to page" right after
#     "Button for initiating a new calculation" (before the blank tab/tab
#     paragraph that follows it).
#  2. Insert a new paragraph "Store variables in array" right before
#     "Write calculated values to output fields".
#  3. Collapse the multi-run "H = (...)" hypergeometric formula paragraph
#     (which had w:proofErr gramStart/gramEnd wrappers splitting it into
#     many runs) into a single run with the full formula text.
#  4. After "For loop to run geometric calculation for success numbers
#     greater than input", insert two new paragraphs:
#       "Inject results graph into website"
#       "Use static image to represent \u201cresults\u201d"
#  5. Add a <w:lastRenderedPageBreak/> marker to the first (tab-only) run
#     of the "Align header to center" paragraph.
#
# Applied from the bottom of the document to the top so paragraph indices
# located ahead of time stay valid as new paragraphs get inserted above
# them.

$d = $word.ActiveDocument
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$CR = [char]13

function Find-ParagraphIndex($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $ptext = $doc.Paragraphs.Item($i).Range.Text.TrimEnd($CR)
        if ($ptext -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# Change 5: lastRenderedPageBreak on "Align header to center" paragraph.
# Content-only edit (paragraph mark / attributes untouched) -- exclude
# the trailing paragraph mark from the InsertXML target range.
# ---------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "`t`tAlign header to center"
if ($idx -eq -1) { throw "Could not find 'Align header to center' paragraph" }
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$xml = $pkgOpen + '<w:p><w:r><w:lastRenderedPageBreak/><w:tab/></w:r><w:r><w:tab/><w:t>Align header to center</w:t></w:r></w:p>' + $pkgClose
$r2.InsertXML($xml)

# ---------------------------------------------------------------------
# Change 4: new paragraphs after "For loop ... greater than input".
# ---------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "`t`tFor loop to run geometric calculation for success numbers greater than input"
if ($idx -eq -1) { throw "Could not find 'For loop ... greater than input' paragraph" }
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$r.InsertParagraphAfter()
$newp = $d.Paragraphs.Item($idx + 1)
$newr = $newp.Range
$xml = $pkgOpen + '<w:p><w:r><w:tab/><w:t>Inject results graph into website</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Use static image to represent &#8220;results&#8221;</w:t></w:r></w:p>' + $pkgClose
$newr.InsertXML($xml)

# ---------------------------------------------------------------------
# Change 3: collapse the "H = (...)" formula runs into one run.
# Content-only edit -- exclude the trailing paragraph mark.
# ---------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "`t`tH = (k!/(x!*(k-x)!)) * ((N-k)!/((n-x)!*((N-k)-(n-x))!)) / (N!/(n!*(N-n)!)"
if ($idx -eq -1) { throw "Could not find 'H = (...)' paragraph" }
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$xml = $pkgOpen + '<w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t>H = (k!/(x!*(k-x)!)) * ((N-k)!/((n-x)!*((N-k)-(n-x))!)) / (N!/(n!*(N-n)!)</w:t></w:r></w:p>' + $pkgClose
$r2.InsertXML($xml)

# ---------------------------------------------------------------------
# Change 2: new paragraph "Store variables in array" before
# "Write calculated values to output fields".
# ---------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "`t`tWrite calculated values to output fields"
if ($idx -eq -1) { throw "Could not find 'Write calculated values...' paragraph" }
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$r.InsertParagraphBefore()
$newp = $d.Paragraphs.Item($idx)
$newr = $newp.Range
$xml = $pkgOpen + '<w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Store variables in array</w:t></w:r></w:p>' + $pkgClose
$newr.InsertXML($xml)

# ---------------------------------------------------------------------
# Change 1: new paragraph "Load results graph to page" before the blank
# tab/tab paragraph that follows "Button for initiating a new
# calculation".
# ---------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "`t`t"
if ($idx -eq -1) { throw "Could not find blank tab/tab paragraph" }
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$r.InsertParagraphBefore()
$newp = $d.Paragraphs.Item($idx)
$newr = $newp.Range
$xml = $pkgOpen + '<w:p><w:r><w:tab/><w:t>Load results graph to page</w:t></w:r></w:p>' + $pkgClose
$newr.InsertXML($xml)

Write-Output "Done"
